$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: becomes a numeric 0, centered horizontally, top-aligned vertically,
# bold font, thin box border all around.
$ws.Range("A1").Value = 0
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1").VerticalAlignment = -4160     # xlTop
$ws.Range("A1").Borders.LineStyle = 1         # xlContinuous
$ws.Range("A1").Borders.Weight = 2            # xlThin

# Rows 2-18: updated text values
$ws.Range("A2").Value = "Dragon Ball Super: Super Heroe"
$ws.Range("A3").Value = "Tren Bala"
$ws.Range("A4").Value = "Minions Nace un Villano"
$ws.Range("A5").Value = "Escalera al Infierno"
$ws.Range("A6").Value = "DC Liga de Super Mascotas"
$ws.Range("A7").Value = "Thor Amor y Trueno"
$ws.Range("A8").Value = "Seventeen World Tour BE THE SUN- Houston LIVE VIEWING"
# A9 "Elvis" is unchanged
$ws.Range("A10").Value = "El Telefono Negro"
$ws.Range("A11").Value = "Alarido"
$ws.Range("A12").Value = "Bestia"
$ws.Range("A13").Value = "Top Gun Maverick"
$ws.Range("A14").Value = "Cambio de Planes"
$ws.Range("A15").Value = "Jurassic World Dominio"
$ws.Range("A16").Value = "2022 Festivales y Muestras"
$ws.Range("A17").Value = "Berta Soy Yo"
$ws.Range("A18").Value = "Indómita Salvaje"

# New rows 19-23
$ws.Range("A19").Value = "Exodo La Ultima Marea"
$ws.Range("A20").Value = "Persiguiendo un Sueño"
$ws.Range("A21").Value = "El Sacrificio"
$ws.Range("A22").Value = "Buena Suerte Leo Grande"
$ws.Range("A23").Value = "El Fotógrafo De Minamata"
